$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the 6 date headers (C2:H2) with 5 homework headers (C2:G2), clear H2
$ws.Range("C2").Value = "ДЗ_1"
$ws.Range("D2").Value = "ДЗ_2"
$ws.Range("E2").Value = "ДЗ_3"
$ws.Range("F2").Value = "ДЗ_4"
$ws.Range("G2").Value = "ДЗ_5"
$ws.Range("H2").Value = $null

# Clear all attendance/grade marks in the data area (rows 4-32, columns C-H)
$ws.Range("C4:H32").Value = $null

# Clear the totals row
$ws.Range("C33:H33").Value = $null

# Update selection to match the author's final selection
$ws.Range("C2:G2").Select()
